$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Antón"
$ws.Range("C3").Value = "Torrón"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Xin"
$ws.Range("C4").Value = "Lú"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Belén"
$ws.Range("C5").Value = "Pastor Iglesias"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Aitor"
$ws.Range("C6").Value = "Menta"

$ws.Columns.Item(3).ColumnWidth = 12.15

$null = $ws.Range("B7").Select()
